$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 0.2062146892655367
    "C2" = 0.5282485875706214
    "J2" = 0.008474576271186441
    "P2" = 0.1581920903954802
    "S2" = 0.09887005649717515
    "B3" = 0.01047120418848168
    "C3" = 0.01047120418848168
    "J3" = 0.03141361256544502
    "P3" = 0.7225130890052356
    "S3" = 0.225130890052356
    "J4" = 0.02040816326530612
    "P4" = 0.6530612244897959
    "S4" = 0.3265306122448979
    "B6" = 0.07471264367816093
    "D6" = 0.02873563218390805
    "F6" = 0.04022988505747126
    "J6" = 0.2241379310344828
    "O6" = 0.02873563218390805
    "Q6" = 0.1264367816091954
    "R6" = 0.07471264367816093
    "S6" = 0.4022988505747127
    "B7" = 0.1586538461538461
    "D7" = 0.01442307692307692
    "F7" = 0.02884615384615385
    "J7" = 0.125
    "O7" = 0.01442307692307692
    "Q7" = 0.1538461538461539
    "R7" = 0.05288461538461538
    "S7" = 0.4519230769230769
    "B8" = 0.1413934426229508
    "D8" = 0.01229508196721311
    "E8" = 0.002049180327868853
    "F8" = 0.06352459016393443
    "J8" = 0.1045081967213115
    "O8" = 0.02049180327868852
    "Q8" = 0.1864754098360656
    "R8" = 0.05737704918032787
    "S8" = 0.4118852459016393
    "B9" = 0.09027777777777778
    "D9" = 0.006944444444444444
    "F9" = 0.03472222222222222
    "J9" = 0.1180555555555556
    "O9" = 0.02083333333333333
    "Q9" = 0.1736111111111111
    "R9" = 0.04166666666666666
    "S9" = 0.5138888888888888
    "B10" = 0.1316025067144136
    "D10" = 0.03043867502238138
    "F10" = 0.05729632945389436
    "J10" = 0.1110116383169203
    "O10" = 0.02327663384064458
    "Q10" = 0.2175470008952551
    "R10" = 0.05640107430617726
    "S10" = 0.3724261414503133
    "G11" = 0.1304347826086956
    "J11" = 0.07023411371237458
    "K11" = 0.1939799331103679
    "L11" = 0.5886287625418061
    "S11" = 0.01672240802675585
    "G12" = 0.7679558011049724
    "J12" = 0.1657458563535912
    "L12" = 0.01104972375690608
    "S12" = 0.05524861878453038
    "F13" = 0.02083333333333333
    "G13" = 0.7291666666666666
    "J13" = 0.2083333333333333
    "S13" = 0.04166666666666666
    "F15" = 0.01015228426395939
    "H15" = 0.1472081218274112
    "I15" = 0.06598984771573604
    "J15" = 0.3045685279187818
    "K15" = 0.07614213197969544
    "M15" = 0.01015228426395939
    "O15" = 0.06091370558375635
    "S15" = 0.3248730964467005
    "F16" = 0.01762114537444934
    "H16" = 0.1629955947136564
    "I16" = 0.06607929515418502
    "J16" = 0.4052863436123348
    "K16" = 0.1101321585903084
    "M16" = 0.03083700440528634
    "N16" = 0.004405286343612335
    "O16" = 0.03083700440528634
    "S16" = 0.1718061674008811
    "F17" = 0.0170316301703163
    "H17" = 0.2214111922141119
    "I17" = 0.07785888077858881
    "J17" = 0.3625304136253041
    "K17" = 0.08029197080291971
    "M17" = 0.0218978102189781
    "O17" = 0.08029197080291971
    "S17" = 0.1386861313868613
    "F18" = 0.03333333333333333
    "H18" = 0.1833333333333333
    "I18" = 0.1
    "J18" = 0.4083333333333333
    "K18" = 0.08333333333333333
    "M18" = 0.01666666666666667
    "O18" = 0.05833333333333333
    "S18" = 0.1166666666666667
    "F19" = 0.01713395638629283
    "H19" = 0.2367601246105919
    "I19" = 0.05841121495327103
    "J19" = 0.3582554517133956
    "K19" = 0.1183800623052959
    "M19" = 0.02336448598130841
    "O19" = 0.05218068535825545
    "S19" = 0.1355140186915888
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

Write-Host "Applied" $updates.Count "cell updates"